# Update stats for 2025-11 (row 24 in the iserv_stats sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B24").Value = 6330
$ws.Range("C24").Value = 1002
$ws.Range("D24").Value = 5928014
$ws.Range("E24").Value = 936.495102685624
$ws.Range("F24").Value = 7.90998977156494
$ws.Range("G24").Value = 3.83419689119171
$ws.Range("H24").Value = 25.57829856212341
